$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 135: copy formatting from the last existing data row (134) ---
$ws.Range("A134:V134").Copy($ws.Range("A135:V135"))

$ws.Range("A135").Value = 134
$ws.Range("B135").Value = "poland"
$ws.Range("C135").Value = "ekstraklasa"
$ws.Range("D135").Value = "2023-2024"
$ws.Range("E135").Value = 45254.75
$ws.Range("F135").Value = "LKS Lodz"
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = "Zaglebie"
$ws.Range("I135").Value = 2
$ws.Range("J135").Value = 3.04
$ws.Range("K135").Value = "16/11/2023 18:12"
$ws.Range("L135").Value = 2.71
$ws.Range("M135").Value = "24/11/2023 17:59"
$ws.Range("N135").Value = 3.4
$ws.Range("O135").Value = "16/11/2023 18:12"
$ws.Range("P135").Value = 3.48
$ws.Range("Q135").Value = "24/11/2023 17:59"
$ws.Range("R135").Value = 2.29
$ws.Range("S135").Value = "16/11/2023 18:12"
$ws.Range("T135").Value = 2.66
$ws.Range("U135").Value = "24/11/2023 17:59"
$ws.Range("V135").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/lks-lodz-zaglebie/Maxb5dUd/"

# --- Row 136: copy formatting from row 134 as well ---
$ws.Range("A134:V134").Copy($ws.Range("A136:V136"))

$ws.Range("A136").Value = 135
$ws.Range("B136").Value = "poland"
$ws.Range("C136").Value = "ekstraklasa"
$ws.Range("D136").Value = "2023-2024"
$ws.Range("E136").Value = 45254.85416666666
$ws.Range("F136").Value = "Jagiellonia"
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = "Piast Gliwice"
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2.78
$ws.Range("K136").Value = "17/11/2023 19:43"
$ws.Range("L136").Value = 2.39
$ws.Range("M136").Value = "24/11/2023 20:29"
$ws.Range("N136").Value = 3.17
$ws.Range("O136").Value = "17/11/2023 19:43"
$ws.Range("P136").Value = 3.12
$ws.Range("Q136").Value = "24/11/2023 20:29"
$ws.Range("R136").Value = 2.6
$ws.Range("S136").Value = "17/11/2023 19:43"
$ws.Range("T136").Value = 3.4
$ws.Range("U136").Value = "24/11/2023 20:29"
$ws.Range("V136").Value = "https://www.betexplorer.com/football/poland/ekstraklasa/jagiellonia-piast-gliwice/tS1PjcF3/"
